# "subi los fsm del RTC"
#
# 1) Bump the fixed ("Update automatically: No") date placeholder text
#    (field type datetimeFigureOut) from 10/09/2016 to 13/09/2016 on the
#    slide master and on every slide layout.
# 2) Nudge four connected FSM shapes on slide 3 (the "Leo minutos" ellipse
#    and the three connectors attached to it).

$p = $ppt.ActivePresentation
$master = $p.SlideMaster
$newDate = "13/09/2016"
$ppPlaceholderDate = 16

for ($i = 1; $i -le $master.Shapes.Count; $i++) {
    $sh = $master.Shapes.Item($i)
    if ($sh.Type -eq 14) {
        if ($sh.PlaceholderFormat.Type -eq $ppPlaceholderDate) {
            $sh.TextFrame.TextRange.Text = $newDate
        }
    }
}

for ($i = 1; $i -le $master.CustomLayouts.Count; $i++) {
    $layout = $master.CustomLayouts.Item($i)
    for ($j = 1; $j -le $layout.Shapes.Count; $j++) {
        $sh = $layout.Shapes.Item($j)
        if ($sh.Type -eq 14) {
            if ($sh.PlaceholderFormat.Type -eq $ppPlaceholderDate) {
                $sh.TextFrame.TextRange.Text = $newDate
            }
        }
    }
}

# --- Slide 3: reposition the FSM shapes around the "Leo minutos" bubble ---
$s3 = $p.Slides.Item(3)

$ellipse = $s3.Shapes.Item("Elipse 14")
$ellipse.Left = 231.44205484409449
$ellipse.Top = 379.6363779527559

$connFromEllipse = $s3.Shapes.Item("Conector recto de flecha 15")
$connFromEllipse.Left = 178.42385826771653
$connFromEllipse.Top = 407.4545746291339

$connIntoEllipse = $s3.Shapes.Item("Conector recto de flecha 18")
$connIntoEllipse.Left = 317.8420563440945
$connIntoEllipse.Top = 295.2
$connIntoEllipse.Width = 1.5761418322834646
$connIntoEllipse.Height = 84.4363779527559

$bentConn = $s3.Shapes.Item("Conector angular 44")
$bentConn.Left = 360.05064392125985
$bentConn.Top = 393.0641732283465
$bentConn.Width = 62.30433090866142
$bentConn.Height = 146.72141732283464
